# Rename worksheet tabs in order to the new sheet names, per diff:
# rerun LU d2c FeatEng for FR cities with new spatial units, and dist models

$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ51831891",
    "summ52238586",
    "summ52543511",
    "summ52838791",
    "summ53132805",
    "summ53431877",
    "summ53716167",
    "summ54033432",
    "summ54331769"
)

for ($i = 0; $i -lt $newNames.Length; $i++) {
    $ws = $wb.Worksheets.Item($i + 1)
    $ws.Name = $newNames[$i]
}
